$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("R2").Value = 12.936
$ws.Range("W2").Value = 0.979
$ws.Range("R3").Value = 0.306
$ws.Range("W3").Value = 0.108
$ws.Range("R4").Value = 2.527
$ws.Range("W4").Value = 0.992
$ws.Range("R5").Value = 4.339
$ws.Range("W5").Value = 0.25
$ws.Range("R6").Value = 645.414
$ws.Range("W6").Value = 0.321
$ws.Range("R7").Value = 29.765
$ws.Range("W7").Value = 0.791
$ws.Range("R8").Value = 33.756
$ws.Range("W8").Value = 0.85
$ws.Range("R9").Value = 12.457
$ws.Range("W9").Value = 0.321
$ws.Range("R10").Value = 1952.056
$ws.Range("W10").Value = 0.99
$ws.Range("R11").Value = 1.911
$ws.Range("W11").Value = 0.991
$ws.Range("R12").Value = 4.928
$ws.Range("W12").Value = 0.928
$ws.Range("R13").Value = 2.046
$ws.Range("W13").Value = 0.976
$ws.Range("R14").Value = 2.05
$ws.Range("W14").Value = 0.982
$ws.Range("R15").Value = 7.336
$ws.Range("R16").Value = 2.15
$ws.Range("R17").Value = 5.124
$ws.Range("W17").Value = 0.97
$ws.Range("R18").Value = 5.18
$ws.Range("W18").Value = 0.972
$ws.Range("R19").Value = 4.839
$ws.Range("W19").Value = 0.97
$ws.Range("R20").Value = 3.86
$ws.Range("W20").Value = 0.887
$ws.Range("R21").Value = 4.417
$ws.Range("W21").Value = 0.871
$ws.Range("R22").Value = 2.652
$ws.Range("R23").Value = 5.561
$ws.Range("R24").Value = 4.975
$ws.Range("W24").Value = 0.996
$ws.Range("R25").Value = 4.995
$ws.Range("R26").Value = 4.641
$ws.Range("R27").Value = 4.87
$ws.Range("W27").Value = 0.968
$ws.Range("R28").Value = 4.36
$ws.Range("R29").Value = 4.986
$ws.Range("W29").Value = 0.958
$ws.Range("R30").Value = 3.933
$ws.Range("R31").Value = 5.05
$ws.Range("R32").Value = 7.043
$ws.Range("R33").Value = 4.244
$ws.Range("W33").Value = 0.937
$ws.Range("R34").Value = 5.257
$ws.Range("W34").Value = 0.979
$ws.Range("R35").Value = 6.31
$ws.Range("W35").Value = 0.99
$ws.Range("R36").Value = 3.518
$ws.Range("R37").Value = 4.32
$ws.Range("W37").Value = 0.98
$ws.Range("R38").Value = 3.495
$ws.Range("R39").Value = 4.627
$ws.Range("W39").Value = 0.929
$ws.Range("R40").Value = 5.485
$ws.Range("W40").Value = 0.89
$ws.Range("R41").Value = 0.255
$ws.Range("W41").Value = 0.991
$ws.Range("R42").Value = 0.907
$ws.Range("W43").Value = 0.667
$ws.Range("R44").Value = 0.437
$ws.Range("W44").Value = 0.976
$ws.Range("R45").Value = 0.105
$ws.Range("W45").Value = 0.996
$ws.Range("W46").Value = 0.462
$ws.Range("R47").Value = 0.95
$ws.Range("R48").Value = 0.013
$ws.Range("R50").Value = 0.007
$ws.Range("R51").Value = 0.012
$ws.Range("R53").Value = 0.086
$ws.Range("W53").Value = 0.99
$ws.Range("R54").Value = 0.007
$ws.Range("R55").Value = 0.001
$ws.Range("R56").Value = 0.913
$ws.Range("R57").Value = 0.019
$ws.Range("R60").Value = 0.001
$ws.Range("R61").Value = 0.016
$ws.Range("R62").Value = 0.009
$ws.Range("R64").Value = 0.766
$ws.Range("W64").Value = 0.893
$ws.Range("R65").Value = 0.867
$ws.Range("W65").Value = 0.995
$ws.Range("R66").Value = 0.465
$ws.Range("R67").Value = 0.295
$ws.Range("W67").Value = 0.864
$ws.Range("R69").Value = 0.338
$ws.Range("W69").Value = 0.433
$ws.Range("R70").Value = 0.641
$ws.Range("R71").Value = 0.108
$ws.Range("W71").Value = 0.353
$ws.Range("R73").Value = 0.108
$ws.Range("W73").Value = 0.26
$ws.Range("R74").Value = 0.262
$ws.Range("W74").Value = 0.99
$ws.Range("R75").Value = 0.1
$ws.Range("W75").Value = 0.261
$ws.Range("R76").Value = 0.731
$ws.Range("W76").Value = 0.99
$ws.Range("R77").Value = 0.062
$ws.Range("W77").Value = 0.896
